$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.701.64'
$ws.Range('E2').Value = '  +1.59%  '
$ws.Range('D3').Value = '3.086.96'
$ws.Range('E3').Value = '  -0.37%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '521.52'
$ws.Range('E5').Value = '  +0.87%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.03'
$ws.Range('E6').Value = '  +0.26%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('E8').Value = '  +0.19%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.35'
$ws.Range('E9').Value = '  +0.01%  '
$ws.Range('E10').Value = '  +0.11%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.384'
$ws.Range('E11').Value = '  +2.42%  '
$ws.Range('D12').Value = '3.619.73'
$ws.Range('E12').Value = '  +0.06%  '
$ws.Range('E13').Value = '  +1.01%  '
$ws.Range('E14').Value = '  +3.21%  '
$ws.Range('E15').Value = '  +0.46%  '
$ws.Range('D16').Value = '58.704.08'
$ws.Range('E16').Value = '  +1.45%  '
$ws.Range('D17').Value = '3.086.79'
$ws.Range('E17').Value = '  -0.37%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.16'
$ws.Range('E18').Value = '  -0.19%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.90'
$ws.Range('E19').Value = '  -1.93%  '
$ws.Range('E20').Value = '  -1.36%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '342.52'
$ws.Range('E21').Value = '  +1.26%  '
$ws.Range('E22').Value = '  -0.22%  '
$ws.Range('E23').Value = '  +0.37%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '65.66'
$ws.Range('E24').Value = '  -0.19%  '
$ws.Range('E25').Value = '  -0.96%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  -0.15%  '
$ws.Range('D27').Value = '0.0₃0921'
$ws.Range('E27').Value = '  -1.65%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.62'
$ws.Range('E28').Value = '  +1.90%  '
$ws.Range('E29').Value = '  +2.02%  '
$ws.Range('E30').Value = '  +1.30%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '21.00'
$ws.Range('E31').Value = '  +0.07%  '
$ws.Range('E32').Value = '  +1.35%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '154.00'
$ws.Range('E33').Value = '  -0.35%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.60'
$ws.Range('E34').Value = '  +0.96%  '
$ws.Range('E35').Value = '  +3.10%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '26.96'
$ws.Range('E36').Value = '  -5.16%  '
$ws.Range('E37').Value = '  +3.55%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0686'
$ws.Range('E38').Value = '  -0.65%  '
$ws.Range('D39').Value = '3.128.11'
$ws.Range('E39').Value = '  -0.35%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.91'
$ws.Range('E40').Value = '  +0.77%  '
$ws.Range('E41').Value = '  -0.64%  '
$ws.Range('E42').Value = '  -0.05%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.664'
$ws.Range('E43').Value = '  -1.62%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.47'
$ws.Range('E44').Value = '  +5.17%  '
$ws.Range('D45').Value = '2.283.98'
$ws.Range('E45').Value = '  -0.19%  '
$ws.Range('E46').Value = '  -0.03%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '20.65'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.958'
$ws.Range('E48').Value = '  +0.50%  '
$ws.Range('E49').Value = '  +1.46%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.747'
$ws.Range('E50').Value = '  +7.70%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '264.47'
$ws.Range('E51').Value = '  +11.00%  '
